$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 22:52:36"
$wsZhCn.Range("H2").Value = "2016-03-18 22:52:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 22:52:39"
$wsDeDe.Range("H2").Value = "2016-03-18 22:53:00"
